$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41: shorten the manual-mode power-cycle instruction text.
$ws.Range("B41").Value = "Cihazı kapatıp aç."

# Rows 42-43 (B42, B43) stay as they are.

# New rows 44-45: insert two new instruction lines.
$ws.Range("B44").Value = "Sürenin dolmasını bekle."
$ws.Range("B45").Value = "Normal şarj(Zmn) iletisini oku."

# Row 46: drop the "Akü hattı kopuk" label in column A and replace column B text.
$ws.Range("A46").ClearContents()
$ws.Range("B46").Value = "Manuel olarak normal şarj moduna geç."

# Row 47: replace text with the new manual-mode message instruction.
$ws.Range("B47").Value = "Normal şarj(Man) iletisini oku."

# New rows 50-51: re-add the battery-fuse note further down the sheet.
$ws.Range("A50").Value = "Akü hattı kopuk"
$ws.Range("B50").Value = "Akü sigorta atık ise test yapılmayacak şekilde ayarlandı."
$ws.Range("B51").Value = "Ancak bazen akü sigorta atıkken de test yapıyor."

# Update the active selection to reflect where the editor ended up.
$ws.Range("B44").Select()
